$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1954.5454
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1875
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1875
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2225

$ws.Range("H76").Value = 4791.727
$ws.Range("I76").Value = 2999
$ws.Range("J76").Value = 5464
$ws.Range("K76").Value = 2999
$ws.Range("L76").Value = 5464
$ws.Range("M76").Value = -2684
$ws.Range("N76").Value = -6094

$ws.Range("H79").Value = 4791.727
$ws.Range("I79").Value = 2999
$ws.Range("J79").Value = 5464
$ws.Range("K79").Value = 2999
$ws.Range("L79").Value = 5464
$ws.Range("M79").Value = -1907
$ws.Range("N79").Value = -7648

$ws.Range("H132").Value = 3127888.2
$ws.Range("I132").Value = 2902.324
$ws.Range("J132").Value = 27780554
$ws.Range("K132").Value = 8706.972
$ws.Range("L132").Value = 83341662
$ws.Range("M132").Value = -6176.972
$ws.Range("N132").Value = -83346722

$ws.Range("H137").Value = 4445039.5
$ws.Range("I137").Value = 496.6111
$ws.Range("J137").Value = 22223210
$ws.Range("K137").Value = 1489.8333
$ws.Range("L137").Value = 66669630
$ws.Range("M137").Value = 1060.1667
$ws.Range("N137").Value = -66674730

$ws.Range("H138").Value = 6668143.5
$ws.Range("I138").Value = 7937751
$ws.Range("J138").Value = 2702.125
$ws.Range("K138").Value = 23813253
$ws.Range("L138").Value = 8106.375
$ws.Range("M138").Value = -23808113
$ws.Range("N138").Value = -18386.375

$ws.Range("H139").Value = 162500
$ws.Range("J139").Value = 162500
$ws.Range("L139").Value = 162500
$ws.Range("N139").Value = -172780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7866.612
$ws.Range("I32").Value = 8256.446
$ws.Range("J32").Value = 6599.65
$ws.Range("K32").Value = 8256.446
$ws.Range("L32").Value = 6599.65
$ws.Range("M32").Value = -7969.446
$ws.Range("N32").Value = -7173.65

$ws.Range("H44").Value = 29916.666
$ws.Range("J44").Value = 29916.666
$ws.Range("L44").Value = 29916.666
$ws.Range("N44").Value = -30892.666

$ws.Range("H55").Value = 30135.8
$ws.Range("J55").Value = 30135.8
$ws.Range("L55").Value = 30135.8
$ws.Range("N55").Value = -30765.8

$ws.Range("H74").Value = 11629452
$ws.Range("I74").Value = 15626502
$ws.Range("K74").Value = 15626502
$ws.Range("M74").Value = -15625628

$ws.Range("H77").Value = 11629452
$ws.Range("I77").Value = 15626502
$ws.Range("K77").Value = 78132510
$ws.Range("M77").Value = -78128142

$ws.Range("H132").Value = 3473472
$ws.Range("I132").Value = 3847270.2
$ws.Range("J132").Value = 2488.2856
$ws.Range("K132").Value = 11541810.6
$ws.Range("L132").Value = 7464.8568
$ws.Range("M132").Value = -11539280.6
$ws.Range("N132").Value = -12524.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4195.5
$ws.Range("I134").Value = 4079.9048
$ws.Range("J134").Value = 5004.6665
$ws.Range("K134").Value = 12239.7144
$ws.Range("L134").Value = 15013.9995
$ws.Range("M134").Value = -9704.714399999999
$ws.Range("N134").Value = -20083.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6412415
$ws.Range("I31").Value = 1880.6487
$ws.Range("J31").Value = 22225066
$ws.Range("K31").Value = 1880.6487
$ws.Range("L31").Value = 22225066
$ws.Range("M31").Value = -1585.6487
$ws.Range("N31").Value = -22225656

$ws.Range("H34").Value = 6412415
$ws.Range("I34").Value = 1880.6487
$ws.Range("J34").Value = 22225066
$ws.Range("K34").Value = 1880.6487
$ws.Range("L34").Value = 22225066
$ws.Range("M34").Value = -1678.6487
$ws.Range("N34").Value = -22225470

$ws.Range("H58").Value = 992.15094
$ws.Range("I58").Value = 507.02127
$ws.Range("K58").Value = 507.02127
$ws.Range("M58").Value = -304.02127

$ws.Range("H86").Value = 3236.3635
$ws.Range("I86").Value = 3014.2856
$ws.Range("J86").Value = 3625
$ws.Range("K86").Value = 3014.2856
$ws.Range("L86").Value = 3625
$ws.Range("M86").Value = -1891.2856
$ws.Range("N86").Value = -5871

$ws.Range("H89").Value = 3236.3635
$ws.Range("I89").Value = 3014.2856
$ws.Range("J89").Value = 3625
$ws.Range("K89").Value = 15071.428
$ws.Range("L89").Value = 18125
$ws.Range("M89").Value = -9455.428
$ws.Range("N89").Value = -29357

$ws.Range("H110").Value = 40368.332
$ws.Range("J110").Value = 40368.332
$ws.Range("L110").Value = 40368.332
$ws.Range("N110").Value = -48548.332

$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -42350

$ws.Range("H120").Value = 39613
$ws.Range("J120").Value = 39613
$ws.Range("L120").Value = 39613
$ws.Range("N120").Value = -46871

$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -32620

$ws.Range("H122").Value = 1825.6666
$ws.Range("I122").Value = 1793.1538
$ws.Range("K122").Value = 5379.4614
$ws.Range("M122").Value = -2929.4614

$ws.Range("H134").Value = 1465.9822
$ws.Range("I134").Value = 1428.8723
$ws.Range("J134").Value = 1659.7778
$ws.Range("K134").Value = 4286.6169
$ws.Range("L134").Value = 4979.3334
$ws.Range("M134").Value = -1751.6169
$ws.Range("N134").Value = -10049.3334

$ws.Range("H136").Value = 992.15094
$ws.Range("I136").Value = 507.02127
$ws.Range("K136").Value = 1521.06381
$ws.Range("M136").Value = 1028.93619

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 852.96155
$ws.Range("I113").Value = 449.85715
$ws.Range("J113").Value = 915.6667
$ws.Range("K113").Value = 1349.57145
$ws.Range("L113").Value = 2747.0001
$ws.Range("M113").Value = 820.4285500000001
$ws.Range("N113").Value = -7087.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 52800
$ws.Range("J138").Value = 52800
$ws.Range("L138").Value = 52800
$ws.Range("N138").Value = -63080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 55780
$ws.Range("J139").Value = 55780
$ws.Range("L139").Value = 55780
$ws.Range("N139").Value = -66060

$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2107.4
$ws.Range("I122").Value = 2419.182
$ws.Range("K122").Value = 7257.545999999999
$ws.Range("M122").Value = -4807.545999999999

$ws.Range("H132").Value = 1633.5366
$ws.Range("I132").Value = 1326.2122
$ws.Range("J132").Value = 2901.25
$ws.Range("K132").Value = 3978.6366
$ws.Range("L132").Value = 8703.75
$ws.Range("M132").Value = -1448.6366
$ws.Range("N132").Value = -13763.75

$ws.Range("H136").Value = 2665.4285
$ws.Range("I136").Value = 1755.8889
$ws.Range("J136").Value = 4302.6
$ws.Range("K136").Value = 5267.6667
$ws.Range("L136").Value = 12907.8
$ws.Range("M136").Value = -2717.6667
$ws.Range("N136").Value = -18007.8
